$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2
Set-TextValue "D2" "42.988.96"
Set-TextValue "E2" "  -5.35%  "

# Row 3
Set-TextValue "D3" "2.218.56"
Set-TextValue "E3" "  -6.55%  "

# Row 4
Set-TextValue "E4" "  -0.07%  "

# Row 5
Set-TextValue "D5" "318.78"
Set-TextValue "E5" "  -3.35%  "

# Row 6
Set-TextValue "D6" "98.13"
Set-TextValue "E6" "  -9.85%  "

# Row 7
Set-TextValue "D7" "0.579"
Set-TextValue "E7" "  -9.16%  "

# Row 8
Set-TextValue "E8" "  -0.10%  "

# Row 9
Set-TextValue "D9" "0.564"
Set-TextValue "E9" "  -8.58%  "

# Row 10
Set-TextValue "D10" "36.48"
Set-TextValue "E10" "  -11.81%  "

# Row 11
Set-TextValue "E11" "  -4.15%  "

# Row 12
Set-TextValue "E12" "  -10.48%  "

# Row 13
Set-TextValue "D13" "7.59"
Set-TextValue "E13" "  -10.93%  "

# Row 14
Set-TextValue "E14" "  -2.85%  "

# Row 15
Set-TextValue "D15" "2.557.01"
Set-TextValue "E15" "  -6.52%  "

# Row 16
Set-TextValue "D16" "0.854"
Set-TextValue "E16" "  -13.30%  "

# Row 17
Set-TextValue "D17" "14.28"
Set-TextValue "E17" "  -7.70%  "

# Row 18
Set-TextValue "D18" "2.220.34"
Set-TextValue "E18" "  -5.80%  "

# Row 19
Set-TextValue "D19" "42.895.40"
Set-TextValue "E19" "  -5.48%  "

# Row 20
Set-TextValue "D20" "13.67"
Set-TextValue "E20" "  -10.92%  "

# Row 21
Set-TextValue "D21" "0.0₃0960"
Set-TextValue "E21" "  -9.96%  "

# Row 22
Set-TextValue "D22" "6.50"
Set-TextValue "E22" "  -11.35%  "

# Row 23
Set-TextValue "D23" "3.22"
Set-TextValue "E23" "  -12.53%  "

# Row 24
Set-TextValue "D24" "65.08"
Set-TextValue "E24" "  -11.27%  "

# Row 25
Set-TextValue "D25" "235.45"
Set-TextValue "E25" "  -10.16%  "

# Row 27
Set-TextValue "E27" "  -0.04%  "

# Row 28
Set-TextValue "D28" "4.04"
Set-TextValue "E28" "  +1.38%  "

# Row 29
Set-TextValue "B29" "Cosmos"
Set-TextValue "C29" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "9.97"
Set-TextValue "E29" "  -12.46%  "

# Row 30
Set-TextValue "B30" "Toncoin"
Set-TextValue "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.22"
Set-TextValue "E30" "  -3.36%  "

# Row 31
Set-TextValue "D31" "6.34"
Set-TextValue "E31" "  -14.78%  "

# Row 32
Set-TextValue "D32" "36.16"
Set-TextValue "E32" "  -3.01%  "

# Row 33
Set-TextValue "D33" "20.18"
Set-TextValue "E33" "  -10.01%  "

# Row 34
Set-TextValue "D34" "0.0859"
Set-TextValue "E34" "  -10.97%  "

# Row 35
Set-TextValue "D35" "154.97"
Set-TextValue "E35" "  -8.06%  "

# Row 36
Set-TextValue "D36" "2.64"
Set-TextValue "E36" "  -6.57%  "

# Row 37
Set-TextValue "E37" "  -2.48%  "

# Row 38
Set-TextValue "D38" "0.121"
Set-TextValue "E38" "  -8.71%  "

# Row 39
Set-TextValue "E39" "  -7.14%  "

# Row 40
Set-TextValue "B40" "Kaspa"
Set-TextValue "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D40" "0.105"
Set-TextValue "E40" "  -11.06%  "

# Row 41
Set-TextValue "B41" "RenderToken"
Set-TextValue "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D41" "4.36"
Set-TextValue "E41" "  -8.77%  "

# Row 42
Set-TextValue "E42" "  -9.54%  "

# Row 43
Set-TextValue "D43" "0.0318"
Set-TextValue "E43" "  -10.66%  "

# Row 44
Set-TextValue "D44" "13.93"
Set-TextValue "E44" "  +6.54%  "

# Row 45
Set-TextValue "E45" "  -0.13%  "

# Row 46
Set-TextValue "D46" "1.733.71"
Set-TextValue "E46" "  -8.15%  "

# Row 47
Set-TextValue "D47" "0.201"
Set-TextValue "E47" "  -12.92%  "

# Row 48
Set-TextValue "D48" "83.59"
Set-TextValue "E48" "  -14.66%  "

# Row 49
Set-TextValue "E49" "  -5.14%  "

# Row 50
Set-TextValue "D50" "5.24"
Set-TextValue "E50" "  -14.78%  "

# Row 51
Set-TextValue "D51" "102.50"
Set-TextValue "E51" "  -9.25%  "
